# Regenerate the localization-status report for the new handoff/handback
# cycle: the source markdown file was renamed from the old GUID
# (d3627a0d-745b-4ce1-9ae3-92e25b28d728) to a new GUID
# (ccf11ad6-82b8-48bf-82a8-0194c3fbf17d), a fresh xliff handoff was
# generated (new hash, new timestamps), and the not-yet-handed-back
# target/handback columns are cleared out pending the new round trip.

$wb = $excel.ActiveWorkbook

$oldGuid = "d3627a0d-745b-4ce1-9ae3-92e25b28d728"
$newGuid = "ccf11ad6-82b8-48bf-82a8-0194c3fbf17d"
$newHash = "03ae10946d8f71a45dfe7efab3526b11daa898f9"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "$newGuid.md"
$overview.Range("B2").Value = "e2e\$newGuid.md"
foreach ($l in $overview.Hyperlinks) {
    $l.TextToDisplay = "e2e\$newGuid.md"
}
$overview.Range("G2").Value = "2016-08-31 17:10:58"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "$newGuid.md"
foreach ($l in $zhcn.Hyperlinks) {
    if ($l.Range.Address() -eq "`$A`$2") {
        $l.TextToDisplay = "$newGuid.md"
    }
}

$zhcn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-08-31 17:10:54"

# Target/handback not produced yet for this round -> clear, and drop
# the stale hyperlink + its hyperlink styling on the old target file.
# (Deleting via a fresh Range(...).Hyperlinks collection wipes every
# hyperlink on the sheet, so walk the sheet collection and only
# delete the item whose Range is I2.)
foreach ($l in $zhcn.Hyperlinks) {
    if ($l.Range.Address() -eq "`$I`$2") {
        $l.Delete()
    }
}
$zhcn.Range("I2").Value = ""
$zhcn.Range("I2").Style = "Normal"
$zhcn.Range("J2").Value = ""
$zhcn.Range("K2").Value = "0001-01-01 00:00:00"

$zhcn.Columns.Item(9).ColumnWidth = 17.8
$zhcn.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "$newGuid.md"
foreach ($l in $dede.Hyperlinks) {
    if ($l.Range.Address() -eq "`$A`$2") {
        $l.TextToDisplay = "$newGuid.md"
    }
}

$dede.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$dede.Range("H2").Value = "2016-08-31 17:10:58"

foreach ($l in $dede.Hyperlinks) {
    if ($l.Range.Address() -eq "`$I`$2") {
        $l.Delete()
    }
}
$dede.Range("I2").Value = ""
$dede.Range("I2").Style = "Normal"
$dede.Range("J2").Value = ""
$dede.Range("K2").Value = "0001-01-01 00:00:00"

$dede.Columns.Item(9).ColumnWidth = 17.8
$dede.Columns.Item(10).ColumnWidth = 20.8
